$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header typo: meal_category -> meal-category
$ws.Range("C1").Value = "meal-category"

# Add new row 3 (mirrors row 2's layout/style)
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A3").Value = "ingredient"
$ws.Range("B3").Value = "https://www.allrecipes.com/recipe/270750/simple-baked-potato/"
$ws.Range("C3").Value = "side-dish"
$ws.Range("D3").Value = "['1 pound baby carrots', '1 1/2 pounds baby red potatoes, scrubbed and halved', '1/4 cup olive oil', '2 teaspoons minced fresh sage', '2 teaspoons minced fresh rosemary', '2 teaspoons minced fresh thyme', '1 teaspoon salt or to taste', '1/2 teaspoon freshly ground black pepper or to taste', '1 tablespoon minced fresh garlic or to taste', '1 onion, cut into eighths', '1 pound zucchini, halved lengthwise and cut into 1-inch pieces', '1/2 pound mushrooms, cleaned and quatered', '1 tablespoon minced fresh parsley, for garnish']"
